$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Supervisor Name value (parallel to Employee Name "P.Mathema" in G4)
$ws.Range("G6").Value = "Prakruti Sinha"
$ws.Range("G6").HorizontalAlignment = -4131

# Supervisor sign-off: initials + date (parallel to employee sign-off in row 25)
$ws.Range("A27").Value = "P.S"
$ws.Range("A27").HorizontalAlignment = -4131

$ws.Range("D27").Value = 41682
$ws.Range("D27").NumberFormat = "m/d/yy"

# Cursor position left where the author clicked last
[void]$ws.Range("I32").Select()
